$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete rows (bottom-up so earlier row numbers stay valid) ---

# Row 13: 004361159 | HFR | 2003.05
$ws.Rows.Item(13).Delete()

# Row 9: 004846293 | LARISSA | 8280.33
$ws.Rows.Item(9).Delete()

# --- Insert rows (still bottom-up relative to remaining inserts) ---

# Insert two new rows before row 5 (004756981 | MATEUS), i.e. right after
# row 4 (005637820 | GUILHERME): 004399832|EULER|40399.8 and 004361159|HFR|33390.27
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "004399832"
$ws.Cells.Item(5, 2).Value = "EULER"
$ws.Cells.Item(5, 3).Value = 40399.8

$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "004361159"
$ws.Cells.Item(6, 2).Value = "HFR"
$ws.Cells.Item(6, 3).Value = 33390.27

# Insert one new row before row 3 (005105172 | VALDIVINO), i.e. right after
# row 2 (004384258 | PAULA): 005529100|DIMITRI|178797.64
$ws.Rows.Item(3).Insert()

$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "005529100"
$ws.Cells.Item(3, 2).Value = "DIMITRI"
$ws.Cells.Item(3, 3).Value = 178797.64
